# Insert one new weekly price record for "Apio" (Terminal Hortofrutícola Agro
# Chillán) as row 91, pushing the existing rows 91-121 down to 92-122.
#
# The new record is for the week of 2021-09-29 (Excel serial 44468):
#   Volumen=160, Precio minimo=9000, Precio maximo=10000,
#   Precio promedio ponderado=9500, Precio $/Kg=1583
# All the other (non-numeric) attributes mirror the record that used to sit
# in row 91 (same market/region/category/variety/quality/unit/origin).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 91:121 down to 92:122, leaving a blank row 91 to populate.
$ws.Rows.Item(91).Insert()

$ws.Range("A91").Value = 7
$ws.Range("B91").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C91").Value = "Ñuble"
$ws.Range("D91").Value = 44468
$ws.Range("E91").Value = 16
$ws.Range("F91").Value = 100112017
$ws.Range("G91").Value = "Apio"
$ws.Range("H91").Value = "Americana (o)"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 160
$ws.Range("K91").Value = 9000
$ws.Range("L91").Value = 10000
$ws.Range("M91").Value = 9500
$ws.Range("N91").Value = "`$/docena de matas"
$ws.Range("O91").Value = "Provincia del Elquí"
$ws.Range("P91").Value = 1583
$ws.Range("Q91").Value = 6
$ws.Range("R91").Value = "Hortaliza"
